$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.411.86'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.566.08'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''1.000'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').Value = '''285.07'
$ws.Range('E6').Value = '  -2.28%  '
$ws.Range('D7').Value = '''0.3634'
$ws.Range('E7').Value = '  -2.55%  '
$ws.Range('D8').Value = '''48.53'
$ws.Range('E8').Value = '  -2.64%  '
$ws.Range('D9').Value = '''0.3343'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').Value = '''1.125'
$ws.Range('E10').Value = '  -1.45%  '
$ws.Range('D11').Value = '''0.07399'
$ws.Range('E11').Value = '  -2.35%  '
$ws.Range('D13').Value = '''20.75'
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('D14').Value = '''5.946'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '''6.895'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '1.566.50'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '''0.00001103'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').Value = '''88.23'
$ws.Range('E18').Value = '  -2.92%  '
$ws.Range('D19').Value = '''0.06703'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').Value = '''6.355'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').Value = '''16.18'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('D24').Value = '22.412.81'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = '''2.390'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').Value = '''2.556'
$ws.Range('E26').Value = '  -4.35%  '
$ws.Range('D27').Value = '''150.06'
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('D28').Value = '''19.34'
$ws.Range('E28').Value = '  -3.69%  '
$ws.Range('D29').Value = '''5.008'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').Value = '''123.71'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').Value = '1.741.48'
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '''1.057'
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('D33').Value = '''2.006'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').Value = '''6.101'
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').Value = '''9.788'
$ws.Range('D36').Value = '''0.08272'
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('D37').Value = '''0.02406'
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('D38').Value = '''0.2228'
$ws.Range('E38').Value = '  -2.40%  '
$ws.Range('D39').Value = '''0.06388'
$ws.Range('E39').Value = '  -2.10%  '
$ws.Range('D40').Value = '''1.289'
$ws.Range('E40').Value = '  -6.42%  '
$ws.Range('D41').Value = '''5.349'
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').Value = '''11.13'
$ws.Range('E42').Value = '  -1.26%  '
$ws.Range('D43').Value = '''0.6089'
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').Value = '''13.77'
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('D46').Value = '''3.756'
$ws.Range('E46').Value = '  -1.54%  '
$ws.Range('D47').Value = '''0.5818'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').Value = '''2.021'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('E49').Value = '  -4.41%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '''0.07201'
$ws.Range('E51').Value = '  -1.68%  '
